$wb = $excel.ActiveWorkbook

# Rename the first sheet from "F-SW-SD-04" to "S-SW-SC-04"
$ws = $wb.Worksheets.Item(1)
$ws.Name = "S-SW-SC-04"

# Move the selection on that sheet from H26 to C3
$ws.Activate()
$ws.Range("C3").Select()
